# Apply the changes described in the commit:
#   "Updated burndown chart to show sprint 2 progress.
#    May need to update sprint 3's planned work as well."
#
# The real (author-entered) edits are two data-entry cells; every other
# difference in the target OOXML (chart number caches, dependent formula
# results, etc.) is a recalculation ripple that Excel derives automatically
# from these two inputs.

$wb = $excel.ActiveWorkbook

# 1) Record Sprint 2's actual realized/completed work on the "PB Burndown"
#    sheet. This is the "show sprint 2 progress" part of the commit -- it
#    feeds the burndown/velocity charts (chart1.xml, chart2.xml) via the
#    ColTopRemainingWork / PBTrend / RealizedSpeed / AverageSpeed* named
#    ranges, so those chart caches update as a side effect of recalculation.
$burndown = $wb.Worksheets.Item("PB Burndown")
$burndown.Range("D29").Value = 21

# 2) Mark the related Product Backlog story (row 23) as "Done" now that the
#    extra sprint-3 planned work has been accounted for.
$productBacklog = $wb.Worksheets.Item("Product Backlog")
$productBacklog.Range("C23").Value = "Done"
